$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -7.808800000000004
$ws.Range("C8").Value = -12.00669999999998
$ws.Range("C10").Value = -13.59189999999999
$ws.Range("C12").Value = -13.27129999999999
$ws.Range("D12").Value = -8.935499999999999
$ws.Range("D15").Value = -8.339799999999991
$ws.Range("D17").Value = -8.408699999999993
$ws.Range("C18").Value = -11.3188
$ws.Range("D26").Value = -7.252800000000003
$ws.Range("D27").Value = -7.480600000000002
$ws.Range("D28").Value = -7.6862
$ws.Range("C37").Value = -13.06620000000001
$ws.Range("D37").Value = -7.262600000000003
$ws.Range("D47").Value = -7.483100000000002
$ws.Range("C55").Value = -13.3826
$ws.Range("D65").Value = -7.672200000000005
$ws.Range("C68").Value = -10.91780000000001
$ws.Range("D73").Value = -8.4392
$ws.Range("C77").Value = -12.69890000000001
$ws.Range("C78").Value = -12.2984
$ws.Range("C81").Value = -14.32519999999999
$ws.Range("C82").Value = -11.9822
$ws.Range("D84").Value = -7.608600000000003
$ws.Range("D85").Value = -8.581099999999998
$ws.Range("D93").Value = -6.748299999999992
$ws.Range("D95").Value = -7.481900000000001
$ws.Range("D98").Value = -7.529800000000002
$ws.Range("D99").Value = -8.120600000000005
$ws.Range("D101").Value = -7.589999999999998
